$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-140 all change from serial date 45192 to 45202
$ws.Range("C2:C140").Value = 45202
